# aggiornamento fino a 02/05
# Appends new daily rows (239-244) to Sheet1, extending the data range
# from A1:D238 to A1:D244, covering 2021-04-27 through 2021-05-02.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: row, date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(239, 44313, 0, 3, 327.5109170305677),
    @(240, 44314, 0, 3, 327.5109170305677),
    @(241, 44315, 1, 3, 327.5109170305677),
    @(242, 44316, 0, 2, 218.3406113537118),
    @(243, 44317, 0, 1, 109.1703056768559),
    @(244, 44318, 0, 1, 109.1703056768559)
)

$lastRow = 238

foreach ($row in $newRows) {
    $r = $row[0]

    # Copy the formatting of column A from the last existing data row so the
    # new date cell keeps the same date style (s="2") used throughout column A.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
